# Update the workbook to add data for 2022-11-13 (extends the "through 11-04"
# running total to "through 11-05").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "as of" date.
$ws.Name = "Through 2022-11-05"

# Update the label in column A for the November row (row 12).
$ws.Range("A12").Value = "November (through 11-05)"

# Update the November row (row 12) values.
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 17
$ws.Range("F12").Value = 7
$ws.Range("G12").Value = 36
$ws.Range("H12").Value = 32
$ws.Range("I12").Value = 14

# Update the Total row (row 13) values.
$ws.Range("B13").Value = 263
$ws.Range("C13").Value = 496
$ws.Range("D13").Value = 730
$ws.Range("E13").Value = 632
$ws.Range("F13").Value = 489
$ws.Range("G13").Value = 1093
$ws.Range("H13").Value = 1473
$ws.Range("I13").Value = 1414
